# Update the AccountCreationData test-data sheet: rename it, refresh the
# sample email addresses, and make it the active/selected sheet (matching
# the "action interface update" tweak from the commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountCreationData")

# Rename the worksheet tab.
$ws.Name = "accountCreationTest"

# Update the email sample data in column A.
$ws.Range("A2").Value = "nht2@gmail.com"
$ws.Range("A3").Value = "qhsd21@gmail.com"
$ws.Range("A4").Value = "jt1@gmail.com"

# Make this sheet the active one, with A4 selected.
$ws.Activate()
$ws.Range("A4").Select()
